# Implemented collision for magic and explosions.
# Changes six matrix cells from "x" (no collision handling) to "f"
# (collision that needs handling), and applies the red "needs handling"
# fill used elsewhere in the sheet for "f" cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellsToUpdate = @("G2", "H2", "B7", "B8", "M8", "H13")

foreach ($addr in $cellsToUpdate) {
    $rng = $ws.Range($addr)
    $rng.Value = "f"
    $rng.Interior.Color = 255
}

# Restore/set the active selection to match the saved workbook state.
$ws.Range("G2").Select()
